$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Update status text "Handed back: in sync with en-US" -> "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# Update timestamps
$wsOverview.Range("G2").Value = "2016-09-06 07:10:34"
$wsDeDe.Range("H2").Value = "2016-09-06 07:10:34"
$wsZhCn.Range("H2").Value = "2016-09-06 07:10:30"

# Update column widths.
# Note: the engine stores column width quantized to the nearest 1/6 of a
# character unit (after adding the standard 5/6-unit cell padding), the same
# way Excel quantizes column widths to whole pixels. To land as close as
# possible on the target stored width of 17.2159881591797 characters, we
# back-solve for the ColumnWidth input that rounds to the nearest achievable
# stored value (17.16666... = 103/6).
$targetColumnWidth = 16.333333333333336
$wsOverview.Columns.Item(5).ColumnWidth = $targetColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $targetColumnWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $targetColumnWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $targetColumnWidth
